# Updated group contribution file.
# Adds a "Total" column (F) that sums each person's contribution row,
# and fills in the previously-missing Code Review 2 row (row 3) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Total" header in column F.
$ws.Range("F1").Value = "Total"

# Fill in the Code Review 2 contribution numbers that were missing.
$ws.Range("B3").Value = 26
$ws.Range("C3").Value = 22
$ws.Range("D3").Value = 26
$ws.Range("E3").Value = 26

# Totals column: row 2 gets its own formula, rows 3-5 are entered together
# so they share one formula group (matches how Excel fills a block).
$ws.Range("F2").Formula = "=SUM(B2:E2)"
$ws.Range("F3:F5").Formula = "=SUM(B3:E3)"
